# Auto-generated edit script applying numeric corrections to multiple sheets
# as described by the target diff (scheduled runner data refresh).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 7412573
$ws.Range("J70").Value = 3728.353
$ws.Range("L70").Value = 11185.059
$ws.Range("N70").Value = -11725.059
$ws.Range("H73").Value = 7412573
$ws.Range("J73").Value = 3728.353
$ws.Range("L73").Value = 11185.059
$ws.Range("N73").Value = -13057.059
$ws.Range("H108").Value = 94975
$ws.Range("J108").Value = 94975
$ws.Range("L108").Value = 94975
$ws.Range("N108").Value = -102655
$ws.Range("H110").Value = 702000000
$ws.Range("J110").Value = 702000000
$ws.Range("L110").Value = 702000000
$ws.Range("N110").Value = -702008180
$ws.Range("H125").Value = 2781.1667
$ws.Range("I125").Value = 629
$ws.Range("K125").Value = 5661
$ws.Range("M125").Value = -3201
$ws.Range("H131").Value = 3631.6365
$ws.Range("I131").Value = 1868.5
$ws.Range("K131").Value = 5605.5
$ws.Range("M131").Value = -565.5
$ws.Range("H132").Value = 3754.718
$ws.Range("I132").Value = 3271.0625
$ws.Range("K132").Value = 9813.1875
$ws.Range("M132").Value = -7283.1875
$ws.Range("H138").Value = 5934.7646
$ws.Range("I138").Value = 1769.9166
$ws.Range("J138").Value = 8206.5
$ws.Range("K138").Value = 5309.7498
$ws.Range("L138").Value = 24619.5
$ws.Range("M138").Value = -169.7497999999996
$ws.Range("N138").Value = -34899.5
$ws.Range("H141").Value = 6667.222
$ws.Range("I141").Value = 4750
$ws.Range("K141").Value = 14250
$ws.Range("M141").Value = -9070

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 23007
$ws.Range("I10").Value = 23007
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 23007
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("M10").Value = -22837
$ws.Range("H45").Value = 84643.03999999999
$ws.Range("I45").Value = 114372.48
$ws.Range("K45").Value = 114372.48
$ws.Range("M45").Value = -113995.48
$ws.Range("H55").Value = 49999.332
$ws.Range("J55").Value = 70000
$ws.Range("L55").Value = 70000
$ws.Range("N55").Value = -70630
$ws.Range("H61").Value = 4827.4194
$ws.Range("I61").Value = 4780.1274
$ws.Range("K61").Value = 4780.1274
$ws.Range("M61").Value = -4568.1274
$ws.Range("H122").Value = 372165.1
$ws.Range("I122").Value = 2867
$ws.Range("K122").Value = 8601
$ws.Range("M122").Value = -6151
$ws.Range("H132").Value = 2892.4375
$ws.Range("I132").Value = 2424.8
$ws.Range("K132").Value = 7274.400000000001
$ws.Range("M132").Value = -4744.400000000001
$ws.Range("H136").Value = 4827.4194
$ws.Range("I136").Value = 4780.1274
$ws.Range("K136").Value = 14340.3822
$ws.Range("M136").Value = -11790.3822

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 21209
$ws.Range("I99").Value = 30669
$ws.Range("K99").Value = 30669
$ws.Range("M99").Value = -29171

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2041.2307
$ws.Range("I16").Value = 2370.7778
$ws.Range("J16").Value = 1299.75
$ws.Range("K16").Value = 2370.7778
$ws.Range("L16").Value = 1299.75
$ws.Range("M16").Value = -2083.7778
$ws.Range("N16").Value = -1873.75
$ws.Range("H31").Value = 2297.2712
$ws.Range("I31").Value = 1464.9333
$ws.Range("J31").Value = 2581.0227
$ws.Range("K31").Value = 1464.9333
$ws.Range("L31").Value = 2581.0227
$ws.Range("M31").Value = -1169.9333
$ws.Range("N31").Value = -3171.0227
$ws.Range("H34").Value = 2297.2712
$ws.Range("I34").Value = 1464.9333
$ws.Range("J34").Value = 2581.0227
$ws.Range("K34").Value = 1464.9333
$ws.Range("L34").Value = 2581.0227
$ws.Range("M34").Value = -1262.9333
$ws.Range("N34").Value = -2985.0227
$ws.Range("H99").Value = 7356506.5
$ws.Range("I99").Value = 15627996
$ws.Range("J99").Value = 4071.4443
$ws.Range("K99").Value = 15627996
$ws.Range("L99").Value = 4071.4443
$ws.Range("M99").Value = -15626498
$ws.Range("N99").Value = -7067.4443
$ws.Range("H112").Value = 49500
$ws.Range("J112").Value = 49500
$ws.Range("L112").Value = 49500
$ws.Range("N112").Value = -52454
$ws.Range("H113").Value = 2041.2307
$ws.Range("I113").Value = 2370.7778
$ws.Range("J113").Value = 1299.75
$ws.Range("K113").Value = 2370.7778
$ws.Range("L113").Value = 1299.75
$ws.Range("M113").Value = -200.7777999999998
$ws.Range("N113").Value = -5639.75
$ws.Range("H126").Value = 7356506.5
$ws.Range("I126").Value = 15627996
$ws.Range("J126").Value = 4071.4443
$ws.Range("K126").Value = 46883988
$ws.Range("L126").Value = 12214.3329
$ws.Range("M126").Value = -46881518
$ws.Range("N126").Value = -17154.3329
$ws.Range("H134").Value = 1607999.1
$ws.Range("I134").Value = 2724357
$ws.Range("J134").Value = 3234.875
$ws.Range("K134").Value = 8173071
$ws.Range("L134").Value = 9704.625
$ws.Range("M134").Value = -8170536
$ws.Range("N134").Value = -14774.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 181.33333
$ws.Range("I2").Value = 128.76923
$ws.Range("J2").Value = 243.45454
$ws.Range("K2").Value = 772.61538
$ws.Range("L2").Value = 1460.72724
$ws.Range("M2").Value = -659.61538
$ws.Range("N2").Value = -1686.72724
$ws.Range("H38").Value = 1636.9445
$ws.Range("I38").Value = 299.4
$ws.Range("J38").Value = 2151.3845
$ws.Range("K38").Value = 898.1999999999999
$ws.Range("L38").Value = 6454.1535
$ws.Range("M38").Value = -551.1999999999999
$ws.Range("N38").Value = -7148.1535
$ws.Range("H92").Value = 927.6923
$ws.Range("I92").Value = 747.1667
$ws.Range("J92").Value = 1082.4286
$ws.Range("K92").Value = 2241.5001
$ws.Range("L92").Value = 3247.2858
$ws.Range("M92").Value = -993.5001000000002
$ws.Range("N92").Value = -5743.2858
$ws.Range("H98").Value = 924.64703
$ws.Range("I98").Value = 1026
$ws.Range("J98").Value = 834.55554
$ws.Range("K98").Value = 3078
$ws.Range("L98").Value = 2503.66662
$ws.Range("M98").Value = -1580
$ws.Range("N98").Value = -5499.66662
$ws.Range("H107").Value = 1726.1613
$ws.Range("J107").Value = 1783.1724
$ws.Range("L107").Value = 5349.5172
$ws.Range("N107").Value = -9189.5172
$ws.Range("H131").Value = 30306664
$ws.Range("I131").Value = 83340220
$ws.Range("J131").Value = 1772.381
$ws.Range("K131").Value = 250020660
$ws.Range("L131").Value = 5317.143
$ws.Range("M131").Value = -250015620
$ws.Range("N131").Value = -15397.143
$ws.Range("H132").Value = 10449120
$ws.Range("I132").Value = 890
$ws.Range("K132").Value = 8010
$ws.Range("M132").Value = -5480
$ws.Range("H140").Value = 9990.096
$ws.Range("I140").Value = 10933
$ws.Range("K140").Value = 32799
$ws.Range("M140").Value = -27619
$ws.Range("H141").Value = 4772.25
$ws.Range("I141").Value = 4772.25
$ws.Range("K141").Value = 14316.75
$ws.Range("M141").Value = -9136.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 22517722
$ws.Range("I11").Value = 28559326
$ws.Range("K11").Value = 28559326
$ws.Range("M11").Value = -28559187
$ws.Range("H12").Value = 1469459.5
$ws.Range("I12").Value = 1645757.6
$ws.Range("K12").Value = 1645757.6
$ws.Range("M12").Value = -1645617.6
$ws.Range("H113").Value = 2701.2856
$ws.Range("I113").Value = 2111
$ws.Range("J113").Value = 2937.4
$ws.Range("K113").Value = 2111
$ws.Range("L113").Value = 2937.4
$ws.Range("M113").Value = 59
$ws.Range("N113").Value = -7277.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8476.483
$ws.Range("J22").Value = 4385.0586
$ws.Range("L22").Value = 4385.0586
$ws.Range("N22").Value = -4975.0586
$ws.Range("H26").Value = 13673
$ws.Range("I26").Value = 20000
$ws.Range("J26").Value = 10509.5
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 10509.5
$ws.Range("M26").Value = -19705
$ws.Range("N26").Value = -11099.5
$ws.Range("H27").Value = 8476.483
$ws.Range("J27").Value = 4385.0586
$ws.Range("L27").Value = 4385.0586
$ws.Range("N27").Value = -4599.0586
$ws.Range("H122").Value = 4291.0938
$ws.Range("I122").Value = 3772.1904
$ws.Range("J122").Value = 5281.727
$ws.Range("K122").Value = 11316.5712
$ws.Range("L122").Value = 15845.181
$ws.Range("M122").Value = -8866.5712
$ws.Range("N122").Value = -20745.181
$ws.Range("H132").Value = 12545.111
$ws.Range("I132").Value = 18614.312
$ws.Range("J132").Value = 3717.182
$ws.Range("K132").Value = 55842.936
$ws.Range("L132").Value = 11151.546
$ws.Range("M132").Value = -53312.936
$ws.Range("N132").Value = -16211.546

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("M12").Value = -1858
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 0
$ws.Range("H13").Value = 3441.3333
$ws.Range("I13").Value = 3441.3333
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 3441.3333
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("M13").Value = -3301.3333
$ws.Range("H81").Value = 27228.25
$ws.Range("I81").Value = 29989.428
$ws.Range("K81").Value = 59978.856
$ws.Range("M81").Value = -58917.856
$ws.Range("H84").Value = 27228.25
$ws.Range("I84").Value = 29989.428
$ws.Range("K84").Value = 299894.28
$ws.Range("M84").Value = -294590.28
$ws.Range("H122").Value = 4676.0386
$ws.Range("I122").Value = 1536.75
$ws.Range("K122").Value = 4610.25
$ws.Range("M122").Value = -2160.25
$ws.Range("H132").Value = 24915.322
$ws.Range("I132").Value = 35208.41
$ws.Range("K132").Value = 105625.23
$ws.Range("M132").Value = -103095.23
